$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "remaining"
$ws.Range("B6").Value = 542.15
